$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.496.18"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.570.90"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  +6.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0591"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0883"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "1.795.55"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "1.571.14"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").Value = "28.497.64"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").Value = "0.0₃0693"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  -5.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "1.392.62"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.534"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.791"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").Value = "1.708.16"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  -4.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0518"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.45%  "
